$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 436, shifting rows 436:496 down to 437:497
$ws.Rows.Item(436).Insert()

# Populate the newly inserted row 436 using values copied down from what is
# now row 437 (the old row 436), except for the columns that actually change.
$ws.Cells.Item(436, 1).Value = 10
$ws.Cells.Item(436, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(436, 3).Value = "La Araucanía"
$ws.Cells.Item(436, 4).Value = 45154
$ws.Cells.Item(436, 4).NumberFormat = $ws.Cells.Item(437, 4).NumberFormat
$ws.Cells.Item(436, 5).Value = 9
$ws.Cells.Item(436, 6).Value = 100112001
$ws.Cells.Item(436, 7).Value = "Berenjena"
$ws.Cells.Item(436, 8).Value = "Sin especificar"
$ws.Cells.Item(436, 9).Value = "Primera"
$ws.Cells.Item(436, 10).Value = 60
$ws.Cells.Item(436, 11).Value = 14000
$ws.Cells.Item(436, 12).Value = 14000
$ws.Cells.Item(436, 13).Value = 14000
$ws.Cells.Item(436, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(436, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(436, 16).Value = 350
$ws.Cells.Item(436, 17).Value = 40
$ws.Cells.Item(436, 18).Value = "Hortaliza"
